$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric "proportion drinking" values (column E, numeric rows) ---
$ws.Range("E2").Value = 0.085811
$ws.Range("E4").Value = 0.085566
$ws.Range("E6").Value = 0.058626
$ws.Range("E8").Value = 0.05577
$ws.Range("E10").Value = 0.049335
$ws.Range("E12").Value = 0.053785
$ws.Range("E14").Value = 0.045293

# --- Text (theta / lambda / proportion drinking) values for "missing" rows ---
$ws.Range("C3").Value = "(0.36)"
$ws.Range("D3").Value = "(0.21)"
$ws.Range("E3").Value = "(0.00001)"

$ws.Range("C5").Value = "(0.36)"
$ws.Range("D5").Value = "(0.24)"
$ws.Range("E5").Value = "(0.00001)"

$ws.Range("C7").Value = "(0.69)"
$ws.Range("D7").Value = "(0.57)"
$ws.Range("E7").Value = "(0.00001)"

$ws.Range("C9").Value = "(0.72)"
$ws.Range("D9").Value = "(0.74)"
$ws.Range("E9").Value = "(0.00001)"

$ws.Range("C11").Value = "(0.93)"
$ws.Range("D11").Value = "(1.06)"
$ws.Range("E11").Value = "(0.00001)"

$ws.Range("C13").Value = "(1.13)"
$ws.Range("D13").Value = "(1.0)"
$ws.Range("E13").Value = "(0.00001)"

$ws.Range("C15").Value = "(1.26)"
$ws.Range("D15").Value = "(0.92)"
$ws.Range("E15").Value = "(0.00001)"
